$wb = $excel.ActiveWorkbook

# Insert the new worksheet right after the first sheet (Workflow_1_TestCases),
# i.e. before the existing "Sequential QA (decision + paral" sheet, and give it
# its final name.
$ws1 = $wb.Worksheets.Item(1)
$new = $wb.Worksheets.Add($null, $ws1)
$new.Name = "Workflow_3_TestCases"

# Header row (row 1) and sample-value row (row 2).
$headers = @(
    "Base Unit of Measure",
    "Volume Unit",
    "External Material Group",
    "Unit of Weight",
    "Plant",
    "MRP Type",
    "Availability Check",
    "MRP Controller",
    "Profit Centre",
    "1st Rem./Exped.",
    "2nd Rem./Exped.",
    "3rd Rem./Exped.",
    "Forecast Model",
    "Storage location",
    "Valuation Type"
)
$values = @(
    "EA -- Each",
    "L -- Liter",
    "00101 -- Finished Products",
    "KG -- Kilogram",
    "0001 -- PLANT 0001",
    "ND -- No planning",
    "01 -- Daily requirements",
    "001 -- JOHN SMITH",
    "PC101 -- Profit Centre PC101",
    "1",
    "2",
    "3",
    "0001 -- M1",
    "0001--M1",
    "N.A -- N.A"
)

# Whole table is formatted as Text, matching the source workbook's style for
# this table (numFmtId 49 / "@").
$new.Range("A1:O3").NumberFormat = "@"

for ($i = 0; $i -lt $headers.Count; $i++) {
    $new.Cells.Item(1, $i + 1).Value = $headers[$i]
}
for ($i = 0; $i -lt $values.Count; $i++) {
    $new.Cells.Item(2, $i + 1).Value = $values[$i]
}

# Row 3 holds plain numbers (0-14); write these before re-applying the text
# format so they stay stored as numeric cells rather than text.
for ($i = 0; $i -lt 15; $i++) {
    $new.Cells.Item(3, $i + 1).Value = $i
}
$new.Range("A3:O3").NumberFormat = "@"

$new.Range("A1:O3").Columns.AutoFit() | Out-Null

$null = $new.Range("H24").Select()
